$wb = $excel.ActiveWorkbook

# --- Sheet "Generator Data" ---
$ws1 = $wb.Worksheets.Item("Generator Data")

$ws1.Range("B2").Value = 70604.9265075
$ws1.Range("C2").Value = 43435.2000634

$ws1.Range("B3").Value = 29661.12962580075
$ws1.Range("C3").Value = 13030.56001902

$ws1.Range("B4").Value = 2966.112962580075
$ws1.Range("C4").Value = 1172.7504017117999

$ws1.Range("B5").Value = 277065.54236800002
$ws1.Range("C5").Value = 23516.467432900001

# --- Sheet "Yearly Fuel Costs" ---
$ws2 = $wb.Worksheets.Item("Yearly Fuel Costs")

$ws2.Range("B2").Value = 150291.0087386785
$ws2.Range("B3").Value = 150291.0010620206
